$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: 29-Sep-2022, 10:00 - 16:00 ---
$ws.Range("B25").Copy($ws.Range("B26"))
$ws.Range("C25").Copy($ws.Range("C26"))
$ws.Range("D25").Copy($ws.Range("D26"))
$ws.Range("B26").Value = 44833
$ws.Range("C26").Value = 0.41666666666666669
$ws.Range("D26").Value = 0.66666666666666663
$ws.Range("E26").Value = "python stepper motor working, made test report, talked with johan brussen about schematic. Cleanup up kicad (libraries where bad) "

# --- Row 27: 30-Sep-2022, "x" "x" "x" ---
$ws.Range("B25").Copy($ws.Range("B27"))
$ws.Range("B27").Value = 44834
$ws.Range("C27").Value = "x"
$ws.Range("D27").Value = "x"
$ws.Range("E27").Value = "x"

# --- Row 28: 3-Oct-2022, 09:45 ---
$ws.Range("B25").Copy($ws.Range("B28"))
$ws.Range("C25").Copy($ws.Range("C28"))
$ws.Range("B28").Value = 44837
$ws.Range("C28").Value = 0.40625
$ws.Range("E28").Value = "watched and learned from videos on python OOP. Tried to create main window, but required data from jeroen. Expirimented with kicad export functionality."

# --- Rows 82-87: jlcpcb part numbers in column D ---
$ws.Range("D82").Value = "C3339"
$ws.Range("D83").Value = "C3339"
$ws.Range("D84").Value = "C160404"
$ws.Range("D85").Value = "C513765"
$ws.Range("D86").Value = "C525005"
$ws.Range("D87").Value = "C2918513"

# --- Update the view: scroll down and select D82:D87 ---
$null = $ws.Range("D82:D87").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
